$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Artn"
$ws.Range("C2").Value = "Ret"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.415594
$ws.Range("H2").Value = 4.246782
$ws.Range("I2").Value = 0.7538443241324221
$ws.Range("J2").Value = 0.7538443241324221
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.21276
$ws.Range("N2").Value = 9.63828
$ws.Range("O2").Value = 0.3952107490920524
$ws.Range("P2").Value = 0.3952107490920524
$ws.Range("Q2").Value = 4.547963779439999
$ws.Range("R2").Value = 40.93167401496
$ws.Range("S2").Value = 0.2979273800391665
$ws.Range("T2").Value = 0.2979273800391665

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Artn"
$ws.Range("C3").Value = "Ret"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.415594
$ws.Range("H3").Value = 4.246782
$ws.Range("I3").Value = 0.7538443241324221
$ws.Range("J3").Value = 0.7538443241324221
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.884996
$ws.Range("N3").Value = 11.654988
$ws.Range("O3").Value = 0.4779044122124365
$ws.Range("P3").Value = 0.4779044122124365
$ws.Range("Q3").Value = 5.499577027623999
$ws.Range("R3").Value = 49.496193248616
$ws.Range("S3").Value = 0.3602655286241866
$ws.Range("T3").Value = 0.3602655286241866

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Artn"
$ws.Range("C4").Value = "Ret"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.415594
$ws.Range("H4").Value = 4.246782
$ws.Range("I4").Value = 0.7538443241324221
$ws.Range("J4").Value = 0.7538443241324221
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.031476333333333
$ws.Range("N4").Value = 3.094429
$ws.Range("O4").Value = 0.1268848386955111
$ws.Range("P4").Value = 0.1268848386955111
$ws.Range("Q4").Value = 1.460151708608667
$ws.Range("R4").Value = 13.141365377478
$ws.Range("S4").Value = 0.09565141546906898
$ws.Range("T4").Value = 0.09565141546906898

$ws.Range("A5").Value = "Inflammatory-Mac"
$ws.Range("B5").Value = "Artn"
$ws.Range("C5").Value = "Ret"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.2356576666666667
$ws.Range("H5").Value = 0.706973
$ws.Range("I5").Value = 0.1254944528268394
$ws.Range("J5").Value = 0.1254944528268394
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.21276
$ws.Range("N5").Value = 9.63828
$ws.Range("O5").Value = 0.3952107490920524
$ws.Range("P5").Value = 0.3952107490920524
$ws.Range("Q5").Value = 0.7571115251599999
$ws.Range("R5").Value = 6.814003726439999
$ws.Range("S5").Value = 0.04959675670859245
$ws.Range("T5").Value = 0.04959675670859245

$ws.Range("A6").Value = "Inflammatory-Mac"
$ws.Range("B6").Value = "Artn"
$ws.Range("C6").Value = "Ret"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.2356576666666667
$ws.Range("H6").Value = 0.706973
$ws.Range("I6").Value = 0.1254944528268394
$ws.Range("J6").Value = 0.1254944528268394
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.884996
$ws.Range("N6").Value = 11.654988
$ws.Range("O6").Value = 0.4779044122124365
$ws.Range("P6").Value = 0.4779044122124365
$ws.Range("Q6").Value = 0.9155290923693332
$ws.Range("R6").Value = 8.239761831324
$ws.Range("S6").Value = 0.05997435271413203
$ws.Range("T6").Value = 0.05997435271413203

$ws.Range("A7").Value = "Inflammatory-Mac"
$ws.Range("B7").Value = "Artn"
$ws.Range("C7").Value = "Ret"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.2356576666666667
$ws.Range("H7").Value = 0.706973
$ws.Range("I7").Value = 0.1254944528268394
$ws.Range("J7").Value = 0.1254944528268394
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.031476333333333
$ws.Range("N7").Value = 3.094429
$ws.Range("O7").Value = 0.1268848386955111
$ws.Range("P7").Value = 0.1268848386955111
$ws.Range("Q7").Value = 0.2430753059352222
$ws.Range("R7").Value = 2.187677753417
$ws.Range("S7").Value = 0.01592334340411495
$ws.Range("T7").Value = 0.01592334340411495

$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Artn"
$ws.Range("C8").Value = "Ret"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.2265816666666667
$ws.Range("H8").Value = 0.679745
$ws.Range("I8").Value = 0.1206612230407385
$ws.Range("J8").Value = 0.1206612230407385
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.21276
$ws.Range("N8").Value = 9.63828
$ws.Range("O8").Value = 0.3952107490920524
$ws.Range("P8").Value = 0.3952107490920524
$ws.Range("Q8").Value = 0.7279525154
$ws.Range("R8").Value = 6.551572638600001
$ws.Range("S8").Value = 0.04768661234429347
$ws.Range("T8").Value = 0.04768661234429347

$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Artn"
$ws.Range("C9").Value = "Ret"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.2265816666666667
$ws.Range("H9").Value = 0.679745
$ws.Range("I9").Value = 0.1206612230407385
$ws.Range("J9").Value = 0.1206612230407385
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.884996
$ws.Range("N9").Value = 11.654988
$ws.Range("O9").Value = 0.4779044122124365
$ws.Range("P9").Value = 0.4779044122124365
$ws.Range("Q9").Value = 0.8802688686733333
$ws.Range("R9").Value = 7.92241981806
$ws.Range("S9").Value = 0.0576645308741178
$ws.Range("T9").Value = 0.0576645308741178

$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Artn"
$ws.Range("C10").Value = "Ret"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.2265816666666667
$ws.Range("H10").Value = 0.679745
$ws.Range("I10").Value = 0.1206612230407385
$ws.Range("J10").Value = 0.1206612230407385
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.031476333333333
$ws.Range("N10").Value = 3.094429
$ws.Range("O10").Value = 0.1268848386955111
$ws.Range("P10").Value = 0.1268848386955111
$ws.Range("Q10").Value = 0.2337136267338889
$ws.Range("R10").Value = 2.103422640605
$ws.Range("S10").Value = 0.01531007982232719
$ws.Range("T10").Value = 0.01531007982232719

Write-Output "done"
